# The source data is a weekly price series. This commit ("Fruta / hortaliza,
# semanal") adds one new week's worth of observations (Primera + Segunda
# quality grades, dated 2022-09-05 / serial 44809) right after the existing
# row for serial 44252 (old row 201). Every row that used to be at 202..275
# shifts down by two to 204..277, which is why the sheet's used range grows
# from A1:R275 to A1:R277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 202 (this pushes the old rows 202-275 down to
# 204-277, matching the dimension change to A1:R277 and every "shifted" cell
# value in the diff).
$ws.Rows.Item(202).Insert()
$ws.Rows.Item(202).Insert()

# New row 202 ("Primera" grade) data
$ws.Cells.Item(202, 1).Value = 7
$ws.Cells.Item(202, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(202, 3).Value = "Ñuble"
$ws.Cells.Item(202, 4).Value = 44809
$ws.Cells.Item(202, 5).Value = 16
$ws.Cells.Item(202, 6).Value = 100112009
$ws.Cells.Item(202, 7).Value = "Acelga"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 240
$ws.Cells.Item(202, 11).Value = 700
$ws.Cells.Item(202, 12).Value = 800
$ws.Cells.Item(202, 13).Value = 750
$ws.Cells.Item(202, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(202, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(202, 16).Value = 750
$ws.Cells.Item(202, 17).Value = 1
$ws.Cells.Item(202, 18).Value = "Hortaliza"

# New row 203 ("Segunda" grade) data
$ws.Cells.Item(203, 1).Value = 7
$ws.Cells.Item(203, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(203, 3).Value = "Ñuble"
$ws.Cells.Item(203, 4).Value = 44809
$ws.Cells.Item(203, 5).Value = 16
$ws.Cells.Item(203, 6).Value = 100112009
$ws.Cells.Item(203, 7).Value = "Acelga"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Segunda"
$ws.Cells.Item(203, 10).Value = 200
$ws.Cells.Item(203, 11).Value = 600
$ws.Cells.Item(203, 12).Value = 600
$ws.Cells.Item(203, 13).Value = 600
$ws.Cells.Item(203, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(203, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(203, 16).Value = 600
$ws.Cells.Item(203, 17).Value = 1
$ws.Cells.Item(203, 18).Value = "Hortaliza"
